# Edit script: Add new SEC enforcement-action rows (19-30), update text for
# rows 3-18 (shared-string index shift after removing the orphaned "Both" value),
# apply custom font styling to H23, and adjust column widths / selection to
# match the author's final saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 'Settlement'
$ws.Range("E3").Value = 'Fraud'
$ws.Range("F3").Value = 'Civil and Criminal'
$ws.Range("G3").Value = 'N/A'
$ws.Range("H3").Value = 'Block Bits Fund'
$ws.Range("I3").Value = 'Bitcoin'
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 'Northern California'

# Row 4
$ws.Range("D4").Value = 'Ongoing'
$ws.Range("E4").Value = 'Fraud'
$ws.Range("F4").Value = 'Civil'
$ws.Range("G4").Value = 'N/A'
$ws.Range("H4").Value = 'MCC International Corp.'
$ws.Range("I4").Value = 'Bitcoin'
$ws.Range("J4").Value = 11300000
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 'Chicago'

# Row 5
$ws.Range("D5").Value = 'Ongoing'
$ws.Range("E5").Value = 'Fraud'
$ws.Range("F5").Value = 'Civil and Criminal'
$ws.Range("G5").Value = 'Ormeus'
$ws.Range("H5").Value = 'Ormeus Coin'
$ws.Range("I5").Value = 'Ethereum'
$ws.Range("J5").Value = 124000000
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 'Southern New York'

# Row 6
$ws.Range("D6").Value = 'Settlement'
$ws.Range("E6").Value = 'Fraud'
$ws.Range("F6").Value = 'Civil'
$ws.Range("G6").Value = 'N/A'
$ws.Range("H6").Value = 'BlockFi'
$ws.Range("I6").Value = 'Bitcoin'
$ws.Range("J6").Value = 100000000
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 'New York'

# Row 7
$ws.Range("D7").Value = 'Ongoing'
$ws.Range("E7").Value = 'Fraud'
$ws.Range("F7").Value = 'Civil'
$ws.Range("G7").Value = 'N/A'
$ws.Range("H7").Value = 'Gold Hawgs'
$ws.Range("I7").Value = 'N/A'
$ws.Range("J7").Value = 400000
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 'Colorado'

# Row 8
$ws.Range("D8").Value = 'Settlement'
$ws.Range("E8").Value = 'Fraud'
$ws.Range("F8").Value = 'Civil'
$ws.Range("G8").Value = 'CMCT'
$ws.Range("H8").Value = 'Crowd Machine'
$ws.Range("I8").Value = 'Ethereum'
$ws.Range("J8").Value = 40700000
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 'Northern California'

# Row 9
$ws.Range("D9").Value = 'Ongoing'
$ws.Range("E9").Value = 'Fraud'
$ws.Range("F9").Value = 'Civil'
$ws.Range("G9").Value = 'DNO'
$ws.Range("H9").Value = 'Denaro'
$ws.Range("I9").Value = 'Ethereum'
$ws.Range("J9").Value = 7000000
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 'New York'

# Row 10
$ws.Range("D10").Value = 'Ongoing'
$ws.Range("E10").Value = 'Fraud'
$ws.Range("F10").Value = 'Civil and Criminal'
$ws.Range("G10").Value = 'N/A'
$ws.Range("H10").Value = 'Social Profimatic and MyMicroProfits.com'
$ws.Range("I10").Value = 'Bitcoin'
$ws.Range("J10").Value = 3600000
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 'Los Angeles'

# Row 11
$ws.Range("D11").Value = 'Ongoing'
$ws.Range("E11").Value = 'Fraud'
$ws.Range("F11").Value = 'Civil'
$ws.Range("G11").Value = 'G-Coins'
$ws.Range("H11").Value = 'GTV Media Group Inc.'
$ws.Range("I11").Value = 'N/A'
$ws.Range("J11").Value = 539000000
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 'New York '

# Row 12
$ws.Range("D12").Value = 'Ongoing'
$ws.Range("E12").Value = 'Unregistered Offering'
$ws.Range("F12").Value = 'Civil'
$ws.Range("G12").Value = 'RvT tokens'
$ws.Range("H12").Value = 'Rivetz Corp.'
$ws.Range("I12").Value = 'Ethereum'
$ws.Range("J12").Value = 18000000
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 'New York'

# Row 13
$ws.Range("D13").Value = 'Ongoing'
$ws.Range("E13").Value = 'Unregistered Offering'
$ws.Range("F13").Value = 'Civil'
$ws.Range("G13").Value = 'N/A'
$ws.Range("H13").Value = 'BitConnect'
$ws.Range("I13").Value = 'Bitcoin'
$ws.Range("J13").Value = 2000000
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 'New York'

# Row 14
$ws.Range("D14").Value = 'Ongoing'
$ws.Range("E14").Value = 'Unregistered Offering'
$ws.Range("F14").Value = 'Civil'
$ws.Range("G14").Value = 'N/A'
$ws.Range("H14").Value = 'BitConnect'
$ws.Range("I14").Value = 'Bitcoin'
$ws.Range("J14").Value = 2000000
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 'New York'

# Row 15
$ws.Range("D15").Value = 'Settlment'
$ws.Range("E15").Value = 'Unregistered Exchange'
$ws.Range("F15").Value = 'Civil'
$ws.Range("G15").Value = 'N/A'
$ws.Range("H15").Value = 'Poloniex'
$ws.Range("I15").Value = 'N/A'
$ws.Range("J15").Value = 10000000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 'Washington, D.C.'

# Row 16
$ws.Range("D16").Value = 'Settlement'
$ws.Range("E16").Value = 'Unregistered Offering'
$ws.Range("F16").Value = 'Civil'
$ws.Range("G16").Value = 'mTokens'
$ws.Range("H16").Value = 'DeFi Money Market'
$ws.Range("I16").Value = 'Ethereum'
$ws.Range("J16").Value = 30000000
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 'New York'

# Row 17
$ws.Range("D17").Value = 'Ongoing'
$ws.Range("E17").Value = 'Unregistered Offering'
$ws.Range("F17").Value = 'Civil'
$ws.Range("G17").Value = 'UULA'
$ws.Range("H17").Value = 'Uulala, Inc.'
$ws.Range("I17").Value = 'N/A'
$ws.Range("J17").Value = 9000000
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 'Los Angeles'

# Row 18
$ws.Range("D18").Value = 'Settlement'
$ws.Range("E18").Value = 'Failed Disclosure'
$ws.Range("F18").Value = 'Civil'
$ws.Range("G18").Value = 'N/A'
$ws.Range("H18").Value = 'Coinschedule'
$ws.Range("I18").Value = 'N/A'
$ws.Range("J18").Value = 197000
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 'Washington, D.C.'

# Row 19
$ws.Range("D19").Value = 'Ongoing'
$ws.Range("E19").Value = 'Unregistered Offering'
$ws.Range("F19").Value = 'Civil'
$ws.Range("G19").Value = 'LOCIcoins'
$ws.Range("H19").Value = 'Loci, Inc.'
$ws.Range("I19").Value = 'Ethereum'
$ws.Range("J19").Value = 7600000
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 'Washington, D.C.'

# Row 20
$ws.Range("D20").Value = 'Settlement'
$ws.Range("E20").Value = 'Unregistered Offering and Fraud'
$ws.Range("F20").Value = 'Civil'
$ws.Range("G20").Value = 'BCT'
$ws.Range("H20").Value = 'CG Blockchain, Inc. and BCT Inc'
$ws.Range("I20").Value = 'N/A'
$ws.Range("J20").Value = 30000000
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 'New York'

# Row 21
$ws.Range("D21").Value = 'Ongoing'
$ws.Range("E21").Value = 'Unregistered Offering and Fraud'
$ws.Range("F21").Value = 'Criminal'
$ws.Range("G21").Value = 'BCT'
$ws.Range("H21").Value = 'CG Blockchain, Inc. and BCT Inc'
$ws.Range("I21").Value = 'N/A'
$ws.Range("J21").Value = 30000000
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 'New Jersey'

# Row 22
$ws.Range("D22").Value = 'Settlement'
$ws.Range("E22").Value = 'Fraud'
$ws.Range("F22").Value = 'Civil'
$ws.Range("G22").Value = 'Apis'
$ws.Range("H22").Value = 'Apis Capital'
$ws.Range("I22").Value = 'Ethereum'
$ws.Range("J22").Value = 20000000
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 'South Carolina'

# Row 23
$ws.Range("D23").Value = 'Ongoing'
$ws.Range("E23").Value = 'Unregistered Offering'
$ws.Range("F23").Value = 'Civil'
$ws.Range("G23").Value = 'LBRY'
$ws.Range("H23").Value = 'LBRY, Inc.'
$ws.Range("I23").Value = 'Ethereum'
$ws.Range("J23").Value = 11000000
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 'New Hampshire'

# Row 24
$ws.Range("D24").Value = 'Ongoing'
$ws.Range("E24").Value = 'Fraud'
$ws.Range("F24").Value = 'Civil'
$ws.Range("G24").Value = 'CTM'
$ws.Range("H24").Value = 'Crypto Traders Management, LLC'
$ws.Range("I24").Value = 'N/A'
$ws.Range("J24").Value = 6900000
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 'Colorado'

# Row 25
$ws.Range("D25").Value = 'Ongoing'
$ws.Range("E25").Value = 'Unregistered Offering'
$ws.Range("F25").Value = 'Civil'
$ws.Range("G25").Value = 'CSD'
$ws.Range("H25").Value = 'Coinseed'
$ws.Range("I25").Value = 'N/A'
$ws.Range("J25").Value = 141410
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 'New York'

# Row 26
$ws.Range("D26").Value = 'Ongoing'
$ws.Range("E26").Value = 'Unregistered Offering and Fraud'
$ws.Range("F26").Value = 'Civil and Criminal'
$ws.Range("G26").Value = 'B2G'
$ws.Range("H26").Value = 'Bitcoiin2Gen'
$ws.Range("I26").Value = 'Ethereum'
$ws.Range("J26").Value = 11000000
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 'New York'

# Row 27
$ws.Range("D27").Value = 'Settlement'
$ws.Range("E27").Value = 'Fraud'
$ws.Range("F27").Value = 'Civil'
$ws.Range("G27").Value = 'N/A'
$ws.Range("H27").Value = 'Wireline, Inc.'
$ws.Range("I27").Value = 'N/A'
$ws.Range("J27").Value = 16000000
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 'Boston'

# Row 28
$ws.Range("D28").Value = 'Settlement'
$ws.Range("E28").Value = 'Unregistered Offering'
$ws.Range("F28").Value = 'Civil'
$ws.Range("G28").Value = 'TNT'
$ws.Range("H28").Value = 'Tierion, Inc.'
$ws.Range("I28").Value = 'N/A'
$ws.Range("J28").Value = 25000000
$ws.Range("K28").Value = 1
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 'Boston'

# Row 29
$ws.Range("D29").Value = 'Ongoing'
$ws.Range("E29").Value = 'Fraud'
$ws.Range("F29").Value = 'Civil'
$ws.Range("G29").Value = 'N/A'
$ws.Range("H29").Value = 'Virgil Capital'
$ws.Range("I29").Value = 'N/A'
$ws.Range("J29").Value = 3500000
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 'San Francisco'

# Row 30
$ws.Range("D30").Value = 'Ongoing'
$ws.Range("E30").Value = 'Unregistered Offering'
$ws.Range("F30").Value = 'Civil'
$ws.Range("G30").Value = 'XRP'
$ws.Range("H30").Value = 'Ripple'
$ws.Range("I30").Value = 'Ripple'
$ws.Range("J30").Value = 1300000000
$ws.Range("K30").Value = 1
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 'New York'

# Special font styling for H23 (LBRY, Inc. cell uses a distinct font color)
$ws.Range("H23").Font.Name = "Times New Roman"
$ws.Range("H23").Font.Size = 12
$ws.Range("H23").Font.Color = 3355443

# Column width adjustments (column E widened, column J widened)
$ws.Columns.Item(5).ColumnWidth = 29.3
$ws.Columns.Item(10).ColumnWidth = 19.65

# Update active selection to match the author's final cursor position
$ws.Range("N30").Select()
